$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AISG Abstract")
Write-Host $ws.Name
